$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: replace the text "good"/"bad" judgement labels with a
#     binary 0/1 "label" column ------------------------------------------------
$ws.Range("B1").Value = "label"

$labels = @{
    2  = 1; 3  = 0; 4  = 0; 5  = 1; 6  = 0; 7  = 1; 8  = 1; 9  = 0; 10 = 0;
    11 = 1; 12 = 0; 13 = 1; 14 = 0; 15 = 0; 16 = 1; 17 = 1; 18 = 1; 19 = 0;
    20 = 0; 21 = 0; 22 = 1; 23 = 0; 24 = 0; 25 = 1; 26 = 1; 27 = 1; 28 = 1;
    29 = 1; 30 = 0; 31 = 0; 32 = 0; 33 = 1; 34 = 1; 35 = 0; 36 = 1
}

foreach ($r in 2..36) {
    $ws.Cells.Item($r, 2).Value = $labels[$r]
}

# the old column used a right-aligned style (s="1"); the new numeric column
# goes back to the default (no explicit cell style)
[void]$ws.Range("B2:B36").ClearFormats()

# --- conditional formatting: the data-bar moved one column to the left ------
$cond = $ws.Range("I24:J25").FormatConditions.Item(1)
$cond.ModifyAppliesToRange($ws.Range("H24:I25"))

# --- selection moved to the last edited cell ---------------------------------
[void]$ws.Range("B36").Select()
